$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert a new row right above "Insurance" (currently row 11) to hold the
# new "cell phone bill" line item. Everything below (Insurance, mastercard,
# earthtreks, Gas, Total Bills, Income, Remainder) shifts down by one, and
# the SUM/"D-minus-D" formulas auto-adjust their ranges.
$ws1.Rows.Item(11).Insert()

# Insert one more blank row before "Gas" (now at row 16) so the gap between
# the earthtreks row and Gas grows from 1 blank row to 2 blank rows.
$ws1.Rows.Item(16).Insert()

# Copy the number formatting (currency style) from the row below onto the
# new Amount cell before putting the value in, so it renders like the rest
# of the "Amount" column.
$ws1.Range("D12").Copy() | Out-Null
$ws1.Range("D11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new "cell phone bill" row: Description, AutoPay?, Amount
$ws1.Range("A11").Value = "cell phone bill"
$ws1.Range("B11").Value = "no"
$ws1.Range("D11").Value = 120

# Matches the author's new cursor position after entering the row.
$ws1.Range("A12").Select()
